# Auto-generated edit script applying the Adamantoise_Profits.xlsx diff
# to the corresponding sheet (named after the in-game crafting class)
# inside the workbook. Values are written cell-by-cell via Range.Value;
# cells that the diff removes entirely are cleared with ClearContents()
# so the <c> element disappears from the saved XML, matching the diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2039.9
$ws.Range("I15").Value = 2039.9
$ws.Range("K15").Value = 6119.700000000001
$ws.Range("M15").Value = -5950.700000000001
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496
$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716
$ws.Range("H92").Value = 21739932
$ws.Range("I92").Value = 26316540
$ws.Range("K92").Value = 26316540
$ws.Range("M92").Value = -26315292
$ws.Range("H96").Value = 1153.75
$ws.Range("I96").Value = 2110
$ws.Range("J96").Value = 197.5
$ws.Range("K96").Value = 6330
$ws.Range("L96").Value = 592.5
$ws.Range("M96").Value = -4957
$ws.Range("N96").Value = -3338.5
$ws.Range("H103").Value = 794.4167
$ws.Range("J103").Value = 872.125
$ws.Range("L103").Value = 2616.375
$ws.Range("N103").Value = -3788.375
$ws.Range("H132").Value = 1835.3158
$ws.Range("I132").Value = 1559.6552
$ws.Range("J132").Value = 2723.5557
$ws.Range("K132").Value = 4678.9656
$ws.Range("L132").Value = 8170.6671
$ws.Range("M132").Value = -2148.9656
$ws.Range("N132").Value = -13230.6671
$ws.Range("H138").Value = 2512.4124
$ws.Range("I138").Value = 1193.1428
$ws.Range("J138").Value = 2734.9397
$ws.Range("K138").Value = 3579.4284
$ws.Range("L138").Value = 8204.819100000001
$ws.Range("M138").Value = 1560.5716
$ws.Range("N138").Value = -18484.8191
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15152514
$ws.Range("I32").Value = 15625890
$ws.Range("J32").Value = 4499
$ws.Range("K32").Value = 15625890
$ws.Range("L32").Value = 4499
$ws.Range("M32").Value = -15625603
$ws.Range("N32").Value = -5073
$ws.Range("H125").Value = 49273
$ws.Range("J125").Value = 49273
$ws.Range("L125").Value = 49273
$ws.Range("N125").Value = -59113
$ws.Range("H132").Value = 3852.3684
$ws.Range("I132").Value = 3824.6875
$ws.Range("K132").Value = 11474.0625
$ws.Range("M132").Value = -8944.0625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1033.1052
$ws.Range("I134").Value = 1033.1052
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3099.3156
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -564.3155999999999
$ws.Range("N134").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5034.325
$ws.Range("I31").Value = 2156.5625
$ws.Range("K31").Value = 2156.5625
$ws.Range("M31").Value = -1861.5625
$ws.Range("H34").Value = 5034.325
$ws.Range("I34").Value = 2156.5625
$ws.Range("K34").Value = 2156.5625
$ws.Range("M34").Value = -1954.5625
$ws.Range("H123").Value = 46124.5
$ws.Range("J123").Value = 46124.5
$ws.Range("L123").Value = 46124.5
$ws.Range("N123").Value = -55924.5
$ws.Range("H132").Value = 6950.9165
$ws.Range("I132").Value = 5548
$ws.Range("K132").Value = 16644
$ws.Range("M132").Value = -14114
$ws.Range("H141").Value = 833992.3
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 833992.3
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 833992.3
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -844352.3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 187803.81
$ws.Range("I11").Value = 200048.6
$ws.Range("J11").Value = 144072.42
$ws.Range("K11").Value = 600145.8
$ws.Range("L11").Value = 432217.26
$ws.Range("M11").Value = -600005.8
$ws.Range("N11").Value = -432497.26
$ws.Range("H34").Value = 787
$ws.Range("I34").Value = 92.14286
$ws.Range("K34").Value = 276.42858
$ws.Range("M34").Value = -192.42858
$ws.Range("H39").Value = 4191.4546
$ws.Range("J39").Value = 4253.875
$ws.Range("L39").Value = 12761.625
$ws.Range("N39").Value = -13349.625
$ws.Range("H52").Value = 1226.75
$ws.Range("J52").Value = 1226.75
$ws.Range("L52").Value = 3680.25
$ws.Range("N52").Value = -4212.25
$ws.Range("H55").Value = 4142.222
$ws.Range("J55").Value = 5018.5713
$ws.Range("L55").Value = 15055.7139
$ws.Range("N55").Value = -15409.7139
$ws.Range("H92").Value = 1707.3334
$ws.Range("I92").Value = 1072.2
$ws.Range("J92").Value = 2501.25
$ws.Range("K92").Value = 3216.6
$ws.Range("L92").Value = 7503.75
$ws.Range("M92").Value = -1968.6
$ws.Range("N92").Value = -9999.75
$ws.Range("H104").Value = 1913.5
$ws.Range("J104").Value = 2500
$ws.Range("L104").Value = 7500
$ws.Range("N104").Value = -12742
$ws.Range("H122").Value = 4999
$ws.Range("J122").Value = 998
$ws.Range("L122").Value = 8982
$ws.Range("N122").Value = -13882
$ws.Range("H131").Value = 1548.25
$ws.Range("I131").Value = 1130.25
$ws.Range("J131").Value = 1966.25
$ws.Range("K131").Value = 3390.75
$ws.Range("L131").Value = 5898.75
$ws.Range("M131").Value = 1649.25
$ws.Range("N131").Value = -15978.75
$ws.Range("H137").Value = 12361.5
$ws.Range("I137").Value = 2019.1666
$ws.Range("J137").Value = 27875
$ws.Range("K137").Value = 6057.4998
$ws.Range("L137").Value = 83625
$ws.Range("M137").Value = -957.4997999999996
$ws.Range("N137").Value = -93825
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3636.1428
$ws.Range("I80").Value = 3463.75
$ws.Range("J80").Value = 3866
$ws.Range("K80").Value = 3463.75
$ws.Range("L80").Value = 3866
$ws.Range("M80").Value = -2465.75
$ws.Range("N80").Value = -5862
$ws.Range("H83").Value = 3636.1428
$ws.Range("I83").Value = 3463.75
$ws.Range("J83").Value = 3866
$ws.Range("K83").Value = 17318.75
$ws.Range("L83").Value = 19330
$ws.Range("M83").Value = -12326.75
$ws.Range("N83").Value = -29314
$ws.Range("H97").Value = 729.6316
$ws.Range("I97").Value = 547.94446
$ws.Range("K97").Value = 547.94446
$ws.Range("M97").Value = -51.94446000000005
$ws.Range("H128").Value = 111329.664
$ws.Range("J128").Value = 111329.664
$ws.Range("L128").Value = 111329.664
$ws.Range("N128").Value = -121289.664
$ws.Range("H130").Value = 121984.5
$ws.Range("J130").Value = 121984.5
$ws.Range("L130").Value = 121984.5
$ws.Range("N130").Value = -132024.5
$ws.Range("H132").Value = 3189.7585
$ws.Range("I132").Value = 2804.375
$ws.Range("J132").Value = 5039.6
$ws.Range("K132").Value = 8413.125
$ws.Range("L132").Value = 15118.8
$ws.Range("M132").Value = -5883.125
$ws.Range("N132").Value = -20178.8
$ws.Range("H134").Value = 83705.25
$ws.Range("J134").Value = 83705.25
$ws.Range("L134").Value = 251115.75
$ws.Range("N134").Value = -256185.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 43479444
$ws.Range("I93").Value = 58824384
$ws.Range("J93").Value = 2109.8333
$ws.Range("K93").Value = 58824384
$ws.Range("L93").Value = 2109.8333
$ws.Range("M93").Value = -58823136
$ws.Range("N93").Value = -4605.8333
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 3798.889
$ws.Range("I132").Value = 2448.8333
$ws.Range("K132").Value = 7346.499899999999
$ws.Range("M132").Value = -4816.499899999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1938.0541
$ws.Range("I132").Value = 1739.9333
$ws.Range("J132").Value = 2787.1428
$ws.Range("K132").Value = 5219.7999
$ws.Range("L132").Value = 8361.428400000001
$ws.Range("M132").Value = -2689.7999
$ws.Range("N132").Value = -13421.4284
